$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2025년 전교생 대상 온라인 SW교육 신청 안내 (유료교육 전액 지원)", "관리자", "2025-03-21", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=38"),
    @("2025년 SW중심대학 주관학과 사업설명회 자료 공유", "관리자", "2025-03-14", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=37"),
    @("2025년도 한국통신학회 동계종합학술발표회 아이디어 경진대회 참가자 SW마일리지 부여", "관리자", "2025-03-14", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=36"),
    @("2025 국민대학교 코딩역량인증제 인증시험 참가자 SW마일리지 부여", "관리자", "2025-03-13", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=35"),
    @("ICCAS 2025 해외교육 프로그램 참가 신청 안내(25.08.25-8.31, 오스트리아/해외교육비 지원)", "관리자", "2025-03-07", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=34"),
    @("2025년 SW중심대학 주관학과 사업설명회 안내(25.03.13(목) 12시~)", "관리자", "2025-03-07", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=33"),
    @("[IITP홍보] 2025 미국 CMU AI 심화교육프로그램 파견 교육생 모집 안내(대학원생-국비 교육비 체제비 지원)", "관리자", "2025-02-10", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=32"),
    @("[소중협] (홍보) 그렙 「2025 프로그래머스 코드챌린지」 개최 안내", "관리자", "2025-01-24", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=31"),
    @("2025년도 SW마에스트로 과정 제16기 연수생 모집 공고 (우수학생 추천)", "관리자", "2025-01-14", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=30")
)

$row = 5
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = "'" + $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}
